# edit.ps1
# Applies odds corrections to Sheet1 of the FlashScore workbook
# as described in the commit diff (row-level value updates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 7).Value = 1.5   # G3
$ws.Cells.Item(3, 8).Value = 3.8   # H3
$ws.Cells.Item(3, 10).Value = 2.1   # J3
$ws.Cells.Item(3, 17).Value = 2.15   # Q3
$ws.Cells.Item(3, 18).Value = 1.67   # R3
$ws.Cells.Item(3, 19).Value = 1.44   # S3
$ws.Cells.Item(3, 20).Value = 2.63   # T3
$ws.Cells.Item(3, 29).Value = 8   # AC3
$ws.Cells.Item(3, 30).Value = 8   # AD3
$ws.Cells.Item(3, 39).Value = 51   # AM3
$ws.Cells.Item(3, 46).Value = 2.63   # AT3

# Row 5
$ws.Cells.Item(5, 7).Value = 1.95   # G5
$ws.Cells.Item(5, 8).Value = 3.4   # H5
$ws.Cells.Item(5, 9).Value = 3.9   # I5
$ws.Cells.Item(5, 10).Value = 2.6   # J5
$ws.Cells.Item(5, 12).Value = 4.33   # L5
$ws.Cells.Item(5, 17).Value = 1.95   # Q5
$ws.Cells.Item(5, 18).Value = 1.9   # R5
$ws.Cells.Item(5, 21).Value = 1.73   # U5
$ws.Cells.Item(5, 22).Value = 2   # V5
$ws.Cells.Item(5, 24).Value = 9.5   # X5
$ws.Cells.Item(5, 25).Value = 9   # Y5
$ws.Cells.Item(5, 26).Value = 17   # Z5
$ws.Cells.Item(5, 34).Value = 12   # AH5
$ws.Cells.Item(5, 35).Value = 21   # AI5
$ws.Cells.Item(5, 36).Value = 13   # AJ5
$ws.Cells.Item(5, 38).Value = 29   # AL5
$ws.Cells.Item(5, 40).Value = 4   # AN5
$ws.Cells.Item(5, 41).Value = 11   # AO5
$ws.Cells.Item(5, 49).Value = 5.5   # AW5
$ws.Cells.Item(5, 50).Value = 21   # AX5
$ws.Cells.Item(5, 52).Value = 67   # AZ5
$ws.Cells.Item(5, 53).Value = 81   # BA5

# Row 7
$ws.Cells.Item(7, 13).Value = 1.1   # M7
$ws.Cells.Item(7, 14).Value = 7   # N7

# Row 10
$ws.Cells.Item(10, 13).Value = 1.03   # M10
$ws.Cells.Item(10, 14).Value = 17   # N10
$ws.Cells.Item(10, 17).Value = 1.57   # Q10
$ws.Cells.Item(10, 18).Value = 2.35   # R10

# Row 13
$ws.Cells.Item(13, 14).Value = 6.8   # N13

# Row 14
$ws.Cells.Item(14, 7).Value = 1.5   # G14
$ws.Cells.Item(14, 9).Value = 6.9   # I14
$ws.Cells.Item(14, 10).Value = 2.02   # J14
$ws.Cells.Item(14, 11).Value = 2.12   # K14
$ws.Cells.Item(14, 14).Value = 6.75   # N14
$ws.Cells.Item(14, 16).Value = 2.75   # P14
$ws.Cells.Item(14, 17).Value = 2   # Q14
$ws.Cells.Item(14, 20).Value = 2.52   # T14
$ws.Cells.Item(14, 21).Value = 2.07   # U14
$ws.Cells.Item(14, 22).Value = 1.6   # V14
$ws.Cells.Item(14, 23).Value = 5.3   # W14
$ws.Cells.Item(14, 24).Value = 6.1   # X14
$ws.Cells.Item(14, 27).Value = 14   # AA14
$ws.Cells.Item(14, 28).Value = 35   # AB14
$ws.Cells.Item(14, 30).Value = 7.2   # AD14
$ws.Cells.Item(14, 31).Value = 20   # AE14
$ws.Cells.Item(14, 32).Value = 120   # AF14
$ws.Cells.Item(14, 36).Value = 22   # AJ14
$ws.Cells.Item(14, 40).Value = 3.15   # AN14
$ws.Cells.Item(14, 41).Value = 7   # AO14
$ws.Cells.Item(14, 43).Value = 22   # AQ14
$ws.Cells.Item(14, 46).Value = 2.5   # AT14
$ws.Cells.Item(14, 47).Value = 8   # AU14
$ws.Cells.Item(14, 48).Value = 90   # AV14
$ws.Cells.Item(14, 49).Value = 7.9   # AW14
